$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = "Volume"
$ws.Range("Q1").Value = "Fragment Size"
$ws.Range("R1").Value = "Read Length"

$ws.Range("P1:R1").Style = $ws.Range("C1").Style

$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 4

$ws.Range("P1").Select()
$ws.Range("P1:R2").Select()
